# Rename the tabs to match their survey code (spaces -> underscores),
# then make "Test_Yearly" (formerly "Test Yearly") the active/selected sheet.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Test Yearly").Name = "Test_Yearly"
$wb.Worksheets.Item("Test Weekly").Name = "Test_Weekly"

$wb.Worksheets.Item("Test_Yearly").Activate()
